$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Update the letter date: "September 19, 2025" -> "September 21, 2025"
# ------------------------------------------------------------------
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2)

# ------------------------------------------------------------------
# 2. Split the mailing address paragraph "979 Story Road, San Jose CA 95122"
#    into two paragraphs:
#       "979 Story Road"
#       "San Jose, CA 95122"
#    (the same text also appears inside a table further down the document
#    and must stay untouched, so we only touch the paragraph that is not
#    inside a table)
# ------------------------------------------------------------------
$addressPara = $null
foreach ($p in $d.Paragraphs) {
    if (($p.Range.Text -like "*979 Story Road, San Jose CA 95122*") -and
        ($p.Range.Information(12) -eq $false)) {
        $addressPara = $p
    }
}

if ($addressPara -ne $null) {
    $addrRange = $addressPara.Range
    $addrRange.Find.Execute(", San Jose CA 95122", $true, $false, $false, $false, $false,
                             $true, 1, $false, "", 2)
    $addressPara.Range.InsertParagraphAfter()
    $addressPara.Next().Range.Text = "San Jose, CA 95122"
}

# ------------------------------------------------------------------
# 3. Remove the empty "No Spacing" paragraph that immediately follows the
#    "Board of Directors" signature line.
# ------------------------------------------------------------------
$boardPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Board of Directors*") {
        $boardPara = $p
    }
}

if ($boardPara -ne $null) {
    $emptyPara = $boardPara.Next()
    if (($emptyPara.Range.Text.Trim().Length -eq 0) -and
        ($emptyPara.Style.NameLocal -eq "No Spacing")) {
        $emptyPara.Range.Delete()
    }
}

Write-Output "edits applied"
